$d = $word.ActiveDocument

function Append-ThreeRuns($paraIndex, $dash, $space, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $e1 = $p.Range.End - 1
    $r1 = $d.Range($e1, $e1)
    $r1.InsertAfter($dash)

    $e2 = $d.Paragraphs.Item($paraIndex).Range.End - 1
    $r2 = $d.Range($e2, $e2)
    $r2.InsertAfter($space)

    $e3 = $d.Paragraphs.Item($paraIndex).Range.End - 1
    $r3 = $d.Range($e3, $e3)
    $r3.InsertAfter($text)
}

# ---------------------------------------------------------------------------
# 1) First "Features to Test" bullet (paragraph 10): replace the dangling
#    run "-login, lists button on nav bar, add " with "-signup, login, buttons"
#    while leaving the "  Features to Test" run untouched (so the run break
#    survives instead of the two runs being coalesced).
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$full10 = $p10.Range.Text
$old10 = "-login, lists button on nav bar, add "
$idx10 = $full10.IndexOf($old10)
$s10 = $p10.Range.Start + $idx10
$e10 = $s10 + $old10.Length
$d.Range($s10, $e10).Delete()
$d.Range($s10, $s10).InsertAfter("-signup, login, buttons")

# ---------------------------------------------------------------------------
# 2) Seven bare "Features to Test" bullets that currently have no trailing
#    run at all: paragraphs 23, 37, 51, 65, 81, 95, 109.
#    Append three new runs: "-", " ", "signup, login, buttons"
# ---------------------------------------------------------------------------
$bareIdxs = @(23, 37, 51, 65, 81, 95, 109)
foreach ($i in $bareIdxs) {
    Append-ThreeRuns $i "-" " " "signup, login, buttons"
}

# ---------------------------------------------------------------------------
# 3) "Features to Test" bullet (paragraph 125) whose trailing run is
#    "-login, " -- delete it and rebuild as three runs: "-", " ",
#    "signup, login, buttons"
# ---------------------------------------------------------------------------
$p125 = $d.Paragraphs.Item(125)
$full125 = $p125.Range.Text
$old125 = "-login, "
$idx125 = $full125.IndexOf($old125)
$s125 = $p125.Range.Start + $idx125
$e125 = $s125 + $old125.Length
$d.Range($s125, $e125).Delete()

$c1 = $s125
$d.Range($c1, $c1).InsertAfter("-")
$c2 = $s125 + 1
$d.Range($c2, $c2).InsertAfter(" ")
$c3 = $s125 + 2
$d.Range($c3, $c3).InsertAfter("signup, login, buttons")

# ---------------------------------------------------------------------------
# 4) Final bare "Features to Test" bullet (paragraph 139). Append the same
#    three runs, then relocate the document's "_GoBack" bookmark from the
#    end of the following "Environments to Test" paragraph to the end of
#    this paragraph (mirrors where Word leaves it after the last edit).
# ---------------------------------------------------------------------------
Append-ThreeRuns 139 "-" " " "signup, login, buttons"

$pEnd = $d.Paragraphs.Item(139).Range.End - 1
$d.Range($pEnd, $pEnd).InsertAfter("Z")
$markRange = $d.Range($pEnd, $pEnd + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$d.Range($pEnd, $pEnd + 1).Text = ""

Write-Output "done"
